# Update of notes 4/23
# Fills in trial data for rows 86-102 (trial_num 85-101) on Sheet1, which
# previously only had the trial_num (column A) and an empty, date-formatted
# column C cell. Also updates the sheet selection to reflect the newly
# populated range.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Each entry: row, pages(B), date serial(C), trial_type(D), checked(E), presence(F), notes(G)
$rows = @(
    @(86,  2, 43195, "M/MLF",            "y", 0, "NA"),
    @(87,  1, 43196, "M/F",              "y", 0, "NA"),
    @(88,  1, 43197, "F/MLF",            "y", 1, "NA"),
    @(89,  1, 43199, "M/MLF",            "y", 0, "NA"),
    @(90,  1, 43200, "M/MLF (assuming)", "y", 0, "redo of previous trial"),
    @(91,  1, 43201, "M/F",              "y", 0, "NA"),
    @(92,  1, 43202, "F/MLF",            "y", 1, "8:13 `"dive but don't think is courtship`""),
    @(93,  1, 43203, "M/MLF",            "y", 0, "NA"),
    @(94,  1, 43204, "M/F",              "y", 1, "NA"),
    @(95,  1, 43205, "F/MLF",            "y", 1, "Side not explicitly specified a few times"),
    @(96,  1, 43206, "M/MLF  ",          "y", 0, "10:36 `"lots of other action`", likely not relevant tho"),
    @(97,  1, 43207, "M/F",              "y", 0, "9:37 check video"),
    @(98,  1, 43208, "F/MLF",            "y", 0, "NA"),
    @(99,  1, 43210, "M/MLF",            "y", 1, "9:06 `"maybe copulate?`""),
    @(100, 2, 43211, "M/F",              "y", 0, "NA"),
    @(101, 2, 43212, "F/MLF",            "y", 0, "NA"),
    @(102, 2, 43213, "M/MLF",            "y", 0, "Second page has nothing on it besides date,trial,etc.")
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 2).Value = $r[1]
    $ws.Cells.Item($rowNum, 3).Value = $r[2]
    $ws.Cells.Item($rowNum, 4).Value = $r[3]
    $ws.Cells.Item($rowNum, 5).Value = $r[4]
    $ws.Cells.Item($rowNum, 6).Value = $r[5]
    $ws.Cells.Item($rowNum, 7).Value = $r[6]
}

# Reflect the new selection/scroll position over the freshly-filled block.
[void]$ws.Range("B86:G102").Select()

Write-Output "updated rows 86-102 on Sheet1"
